$d = $word.ActiveDocument

# Replace every occurrence of $old with $new, anywhere in the document
# body. We wrap the matched span in a throw-away bookmark and set the
# text through the bookmark's Range rather than through a bare
# Document.Range(...): doing so updates only the run(s) that actually
# contain the text and leaves neighbouring (e.g. empty formatting)
# runs in the same paragraph completely untouched.
function Replace-AllOccurrences($old, $new) {
    $bmCounter = 0
    while ($true) {
        $full = $d.Content.Text
        $idx = $full.IndexOf($old)
        if ($idx -lt 0) {
            break
        }
        $bmCounter = $bmCounter + 1
        $target = $d.Range($idx, $idx + $old.Length)
        $bmName = "tmpRepl" + $bmCounter
        $d.Bookmarks.Add($bmName, $target) | Out-Null
        $bm = $d.Bookmarks.Item($bmName)
        $bm.Range.Text = $new
        $d.Bookmarks.Item($bmName).Delete()
    }
}

# Title (appears twice: main heading + bolded run near the end)
Replace-AllOccurrences "Play Dragon Spin Pick n Mix Free - Review of Features & Payouts" "Play Dragon Spin Pick n Mix for Free"

# "What we like" bullet list
Replace-AllOccurrences "Multiple bonus features" "Exciting and unique bonus features"
Replace-AllOccurrences "Excellent graphics and music" "High-quality graphics and immersive music"
Replace-AllOccurrences "Accessible to all with varying bet sizes" "High volatility for experienced players"
Replace-AllOccurrences "High volatility for high-stakes gambling" "Multiple betting options and potential for big payouts"

# "What we don't like" bullet list
Replace-AllOccurrences "Fixed number of game lines" "Limited betting range (€0.10 to €100)"
Replace-AllOccurrences "Gameplay may be too risky for some" "Can be overwhelming for new players"

# Meta description (italic run)
Replace-AllOccurrences "Read a review of Dragon Spin Pick n Mix, an online slot game with varying bet sizes and high volatility. Play for free and activate multiple bonuses." "Read our review of Dragon Spin Pick n Mix and play for free. Discover exciting bonuses and high-quality graphics."
